# Add dil_4..dil_8 columns (K:O) to the "lung" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lung")

# Header row (row 1) - new dilution columns
$ws.Range("K1").Value = "dil_4"
$ws.Range("L1").Value = "dil_5"
$ws.Range("M1").Value = "dil_6"
$ws.Range("N1").Value = "dil_7"
$ws.Range("O1").Value = "dil_8"

# Data rows 2-31: fill new columns with 0, centered like the other dilution columns
$dataRange = $ws.Range("K2:O31")
$dataRange.Value = 0
$dataRange.HorizontalAlignment = -4108

# Update selection to match the target state
[void]$ws.Range("K44").Select()
